# Updates the cryptocurrency price/volume table to reflect the latest
# scraped values from coinranking.com (GitHub Actions scheduled refresh).
#
# Column D (Price) values that look like plain numbers are written with a
# leading apostrophe so Excel keeps them as text (preserving formats such
# as trailing zeros, e.g. "1.00", "0.200") instead of silently coercing
# them to numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.067.39'
$ws.Range("E2").Value = '  -2.40%  '
$ws.Range("D3").Value = '2.234.63'
$ws.Range("E3").Value = '  -3.16%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = "'" + '245.48'
$ws.Range("E5").Value = '  -2.96%  '
$ws.Range("D6").Value = "'" + '0.621'
$ws.Range("E6").Value = '  -3.51%  '
$ws.Range("D7").Value = "'" + '76.13'
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  -5.61%  '
$ws.Range("D10").Value = "'" + '41.22'
$ws.Range("E10").Value = '  +3.88%  '
$ws.Range("D11").Value = "'" + '0.0942'
$ws.Range("E11").Value = '  -4.13%  '
$ws.Range("D12").Value = "'" + '6.99'
$ws.Range("E12").Value = '  -10.09%  '
$ws.Range("E13").Value = '  -2.90%  '
$ws.Range("D14").Value = '2.568.76'
$ws.Range("E14").Value = '  -3.06%  '
$ws.Range("D15").Value = "'" + '14.60'
$ws.Range("E15").Value = '  -6.05%  '
$ws.Range("E16").Value = '  -3.26%  '
$ws.Range("D17").Value = '2.237.28'
$ws.Range("E17").Value = '  -2.47%  '
$ws.Range("D18").Value = '41.929.30'
$ws.Range("E18").Value = '  -2.64%  '
$ws.Range("D19").Value = '0.0₃0977'
$ws.Range("E19").Value = '  -3.93%  '
$ws.Range("D20").Value = "'" + '71.49'
$ws.Range("E20").Value = '  -2.18%  '
$ws.Range("D21").Value = "'" + '6.05'
$ws.Range("E21").Value = '  -4.10%  '
$ws.Range("D22").Value = "'" + '2.27'
$ws.Range("E22").Value = '  +0.37%  '
$ws.Range("D23").Value = "'" + '229.98'
$ws.Range("E23").Value = '  -3.77%  '
$ws.Range("D24").Value = "'" + '1.00'
$ws.Range("E24").Value = '  -0.02%  '
$ws.Range("D25").Value = "'" + '3.68'
$ws.Range("E25").Value = '  -5.84%  '
$ws.Range("D26").Value = "'" + '11.15'
$ws.Range("E26").Value = '  -4.36%  '
$ws.Range("E27").Value = '  -6.35%  '
$ws.Range("D28").Value = "'" + '7.34'
$ws.Range("E28").Value = '  +14.65%  '
$ws.Range("E29").Value = '  -2.01%  '
$ws.Range("D30").Value = "'" + '168.46'
$ws.Range("E30").Value = '  +0.31%  '
$ws.Range("D31").Value = "'" + '20.48'
$ws.Range("E31").Value = '  -3.84%  '
$ws.Range("D32").Value = "'" + '0.0823'
$ws.Range("E32").Value = '  -2.68%  '
$ws.Range("D33").Value = "'" + '32.08'
$ws.Range("E33").Value = '  +4.40%  '
$ws.Range("E34").Value = '  -7.64%  '
$ws.Range("E35").Value = '  -1.77%  '
$ws.Range("D36").Value = "'" + '4.43'
$ws.Range("E36").Value = '  -4.18%  '
$ws.Range("D37").Value = "'" + '4.91'
$ws.Range("E37").Value = '  +1.16%  '
$ws.Range("D38").Value = "'" + '0.0298'
$ws.Range("E38").Value = '  -4.87%  '
$ws.Range("D39").Value = "'" + '13.91'
$ws.Range("E39").Value = '  -0.13%  '
$ws.Range("E40").Value = '  -1.55%  '
$ws.Range("E41").Value = '  -8.52%  '
$ws.Range("D42").Value = "'" + '112.31'
$ws.Range("E42").Value = '  +6.88%  '
$ws.Range("D43").Value = "'" + '0.200'
$ws.Range("E43").Value = '  -8.62%  '
$ws.Range("D44").Value = "'" + '60.26'
$ws.Range("E44").Value = '  -4.05%  '
$ws.Range("E45").Value = '  -6.53%  '
$ws.Range("D46").Value = "'" + '0.0988'
$ws.Range("E46").Value = '  -4.63%  '
$ws.Range("E47").Value = '  -0.49%  '
$ws.Range("D48").Value = "'" + '1.12'
$ws.Range("E48").Value = '  -5.31%  '
$ws.Range("E49").Value = '  -3.08%  '
$ws.Range("D50").Value = "'" + '0.429'
$ws.Range("E50").Value = '  +12.75%  '
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").Value = "'" + '2.24'
$ws.Range("E51").Value = '  -3.62%  '
